$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Cells.Item(28, 8).Value = 3088.9473
$ws.Cells.Item(28, 9).Value = 2850.0667
$ws.Cells.Item(28, 10).Value = 3984.75
$ws.Cells.Item(28, 11).Value = 2850.0667
$ws.Cells.Item(28, 12).Value = 3984.75
$ws.Cells.Item(28, 13).Value = -2365.0667
$ws.Cells.Item(28, 14).Value = -4954.75
# row 38
$ws.Cells.Item(38, 8).Value = 1979.24
$ws.Cells.Item(38, 9).Value = 104.888885
$ws.Cells.Item(38, 10).Value = 3033.5625
$ws.Cells.Item(38, 11).Value = 314.666655
$ws.Cells.Item(38, 12).Value = 9100.6875
$ws.Cells.Item(38, 13).Value = 57.33334500000001
$ws.Cells.Item(38, 14).Value = -9844.6875
# row 40
$ws.Cells.Item(40, 8).Value = 3000
$ws.Cells.Item(40, 10).Value = 3000
$ws.Cells.Item(40, 12).Value = 3000
$ws.Cells.Item(40, 14).Value = -3350
# row 98
$ws.Cells.Item(98, 8).Value = 8729.429
$ws.Cells.Item(98, 10).Value = 1506
$ws.Cells.Item(98, 12).Value = 1506
$ws.Cells.Item(98, 14).Value = -4502
# row 112
$ws.Cells.Item(112, 8).Value = 2647.889
$ws.Cells.Item(112, 10).Value = 3204.6428
$ws.Cells.Item(112, 12).Value = 9613.928400000001
$ws.Cells.Item(112, 14).Value = -11829.9284
# row 122
$ws.Cells.Item(122, 8).Value = 8729.429
$ws.Cells.Item(122, 10).Value = 1506
$ws.Cells.Item(122, 12).Value = 4518
$ws.Cells.Item(122, 14).Value = -9418
# row 129
$ws.Cells.Item(129, 8).Value = 778.0244
$ws.Cells.Item(129, 10).Value = 860.94116
$ws.Cells.Item(129, 12).Value = 2582.82348
$ws.Cells.Item(129, 14).Value = -12582.82348
# row 138
$ws.Cells.Item(138, 8).Value = 1749.34
$ws.Cells.Item(138, 9).Value = 1027.6957
$ws.Cells.Item(138, 10).Value = 1964.8961
$ws.Cells.Item(138, 11).Value = 3083.0871
$ws.Cells.Item(138, 12).Value = 5894.6883
$ws.Cells.Item(138, 13).Value = 2056.9129
$ws.Cells.Item(138, 14).Value = -16174.6883

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Cells.Item(32, 8).Value = 7452.7
$ws.Cells.Item(32, 9).Value = 5701.0815
$ws.Cells.Item(32, 11).Value = 5701.0815
$ws.Cells.Item(32, 13).Value = -5414.0815
# row 45
$ws.Cells.Item(45, 8).Value = 1311.8182
$ws.Cells.Item(45, 9).Value = 1305.7142
$ws.Cells.Item(45, 10).Value = 1322.5
$ws.Cells.Item(45, 11).Value = 1305.7142
$ws.Cells.Item(45, 12).Value = 1322.5
$ws.Cells.Item(45, 13).Value = -928.7141999999999
$ws.Cells.Item(45, 14).Value = -2076.5
# row 102
$ws.Cells.Item(102, 8).Value = 11906371
$ws.Cells.Item(102, 9).Value = 15152790
$ws.Cells.Item(102, 11).Value = 15152790
$ws.Cells.Item(102, 13).Value = -15151168
# row 122
$ws.Cells.Item(122, 8).Value = 2092.682
$ws.Cells.Item(122, 9).Value = 1858.4736
$ws.Cells.Item(122, 10).Value = 3576
$ws.Cells.Item(122, 11).Value = 5575.4208
$ws.Cells.Item(122, 12).Value = 10728
$ws.Cells.Item(122, 13).Value = -3125.4208
$ws.Cells.Item(122, 14).Value = -15628
# row 132
$ws.Cells.Item(132, 8).Value = 2392.4792
$ws.Cells.Item(132, 9).Value = 1617.9688
$ws.Cells.Item(132, 10).Value = 3941.5
$ws.Cells.Item(132, 11).Value = 4853.9064
$ws.Cells.Item(132, 12).Value = 11824.5
$ws.Cells.Item(132, 13).Value = -2323.9064
$ws.Cells.Item(132, 14).Value = -16884.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Cells.Item(94, 8).Value = 9259663
$ws.Cells.Item(94, 9).Value = 9615727
$ws.Cells.Item(94, 10).Value = 2000
$ws.Cells.Item(94, 11).Value = 9615727
$ws.Cells.Item(94, 12).Value = 2000
$ws.Cells.Item(94, 13).Value = -9615276
$ws.Cells.Item(94, 14).Value = -2902
# row 99
$ws.Cells.Item(99, 8).Value = 62501012
$ws.Cells.Item(99, 9).Value = 71429510
$ws.Cells.Item(99, 10).Value = 1490
$ws.Cells.Item(99, 11).Value = 71429510
$ws.Cells.Item(99, 12).Value = 1490
$ws.Cells.Item(99, 13).Value = -71428012
$ws.Cells.Item(99, 14).Value = -4486
# row 105
$ws.Cells.Item(105, 8).Value = 76924490
$ws.Cells.Item(105, 9).Value = 100001110
$ws.Cells.Item(105, 11).Value = 100001110
$ws.Cells.Item(105, 13).Value = -99999363
# row 134
$ws.Cells.Item(134, 8).Value = 1213.2106
$ws.Cells.Item(134, 9).Value = 939.4286
$ws.Cells.Item(134, 11).Value = 2818.2858
$ws.Cells.Item(134, 13).Value = -283.2857999999997

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Cells.Item(31, 8).Value = 1157.4445
$ws.Cells.Item(31, 9).Value = 1105.9
$ws.Cells.Item(31, 10).Value = 1415.1666
$ws.Cells.Item(31, 11).Value = 1105.9
$ws.Cells.Item(31, 12).Value = 1415.1666
$ws.Cells.Item(31, 13).Value = -810.9000000000001
$ws.Cells.Item(31, 14).Value = -2005.1666
# row 34
$ws.Cells.Item(34, 8).Value = 1157.4445
$ws.Cells.Item(34, 9).Value = 1105.9
$ws.Cells.Item(34, 10).Value = 1415.1666
$ws.Cells.Item(34, 11).Value = 1105.9
$ws.Cells.Item(34, 12).Value = 1415.1666
$ws.Cells.Item(34, 13).Value = -903.9000000000001
$ws.Cells.Item(34, 14).Value = -1819.1666
# row 132
$ws.Cells.Item(132, 8).Value = 2040.5217
$ws.Cells.Item(132, 9).Value = 1403.8572
$ws.Cells.Item(132, 10).Value = 3030.889
$ws.Cells.Item(132, 11).Value = 4211.571599999999
$ws.Cells.Item(132, 12).Value = 9092.667000000001
$ws.Cells.Item(132, 13).Value = -1681.571599999999
$ws.Cells.Item(132, 14).Value = -14152.667
# row 140
$ws.Cells.Item(140, 8).Value = 38919.8
$ws.Cells.Item(140, 10).Value = 38919.8
$ws.Cells.Item(140, 12).Value = 38919.8
$ws.Cells.Item(140, 14).Value = -49279.8

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 74
$ws.Cells.Item(74, 8).Value = 5599.8
$ws.Cells.Item(74, 10).Value = 5599.8
$ws.Cells.Item(74, 12).Value = 16799.4
$ws.Cells.Item(74, 14).Value = -18921.4
# row 77
$ws.Cells.Item(77, 8).Value = 5599.8
$ws.Cells.Item(77, 10).Value = 5599.8
$ws.Cells.Item(77, 12).Value = 50398.2
$ws.Cells.Item(77, 14).Value = -61006.2
# row 120
$ws.Cells.Item(120, 8).Value = 6032.2
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 6032.2
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 18096.6
$ws.Cells.Item(120, 13).Value = $null
$ws.Cells.Item(120, 14).Value = -27772.6
# row 121
$ws.Cells.Item(121, 8).Value = 830.6923
$ws.Cells.Item(121, 9).Value = 359.8
$ws.Cells.Item(121, 10).Value = 1125
$ws.Cells.Item(121, 11).Value = 1079.4
$ws.Cells.Item(121, 12).Value = 3375
$ws.Cells.Item(121, 13).Value = 230.5999999999999
$ws.Cells.Item(121, 14).Value = -5995
# row 124
$ws.Cells.Item(124, 8).Value = 1655.5555
$ws.Cells.Item(124, 9).Value = 450
$ws.Cells.Item(124, 10).Value = 2000
$ws.Cells.Item(124, 11).Value = 1350
$ws.Cells.Item(124, 12).Value = 6000
$ws.Cells.Item(124, 13).Value = 3560
$ws.Cells.Item(124, 14).Value = -15820
# row 125
$ws.Cells.Item(125, 8).Value = 5021
$ws.Cells.Item(125, 9).Value = 2030
$ws.Cells.Item(125, 10).Value = 6516.5
$ws.Cells.Item(125, 11).Value = 6090
$ws.Cells.Item(125, 12).Value = 19549.5
$ws.Cells.Item(125, 13).Value = -1170
$ws.Cells.Item(125, 14).Value = -29389.5
# row 129
$ws.Cells.Item(129, 8).Value = 23149058
$ws.Cells.Item(129, 10).Value = 5953439.5
$ws.Cells.Item(129, 12).Value = 17860318.5
$ws.Cells.Item(129, 14).Value = -17870318.5
# row 131
$ws.Cells.Item(131, 8).Value = 23258996
$ws.Cells.Item(131, 10).Value = 4149.5625
$ws.Cells.Item(131, 12).Value = 12448.6875
$ws.Cells.Item(131, 14).Value = -22528.6875
# row 132
$ws.Cells.Item(132, 8).Value = 1097.6428
$ws.Cells.Item(132, 10).Value = 1375.8
$ws.Cells.Item(132, 12).Value = 12382.2
$ws.Cells.Item(132, 14).Value = -17442.2
# row 134
$ws.Cells.Item(134, 8).Value = 2778.8845
$ws.Cells.Item(134, 9).Value = 2845.9
$ws.Cells.Item(134, 10).Value = 2737
$ws.Cells.Item(134, 11).Value = 8537.700000000001
$ws.Cells.Item(134, 12).Value = 8211
$ws.Cells.Item(134, 13).Value = -3467.700000000001
$ws.Cells.Item(134, 14).Value = -18351
# row 136
$ws.Cells.Item(136, 8).Value = 2609.375
$ws.Cells.Item(136, 9).Value = 2213.111
$ws.Cells.Item(136, 10).Value = 3118.8572
$ws.Cells.Item(136, 11).Value = 6639.333
$ws.Cells.Item(136, 12).Value = 9356.571599999999
$ws.Cells.Item(136, 13).Value = -1539.333
$ws.Cells.Item(136, 14).Value = -19556.5716
# row 137
$ws.Cells.Item(137, 8).Value = 32611652
$ws.Cells.Item(137, 9).Value = 93751730
$ws.Cells.Item(137, 10).Value = 3611.0667
$ws.Cells.Item(137, 11).Value = 281255190
$ws.Cells.Item(137, 12).Value = 10833.2001
$ws.Cells.Item(137, 13).Value = -281250090
$ws.Cells.Item(137, 14).Value = -21033.2001
# row 139
$ws.Cells.Item(139, 8).Value = 1674.08
$ws.Cells.Item(139, 9).Value = 1470.15
$ws.Cells.Item(139, 10).Value = 1810.0333
$ws.Cells.Item(139, 11).Value = 4410.450000000001
$ws.Cells.Item(139, 12).Value = 5430.0999
$ws.Cells.Item(139, 13).Value = 729.5499999999993
$ws.Cells.Item(139, 14).Value = -15710.0999
# row 140
$ws.Cells.Item(140, 8).Value = 23226.184
$ws.Cells.Item(140, 9).Value = 49576.668
$ws.Cells.Item(140, 11).Value = 148730.004
$ws.Cells.Item(140, 13).Value = -143550.004
# row 141
$ws.Cells.Item(141, 8).Value = 71431736
$ws.Cells.Item(141, 9).Value = 83335190
$ws.Cells.Item(141, 10).Value = 11016.5
$ws.Cells.Item(141, 11).Value = 250005570
$ws.Cells.Item(141, 12).Value = 33049.5
$ws.Cells.Item(141, 13).Value = -250000390
$ws.Cells.Item(141, 14).Value = -43409.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 113
$ws.Cells.Item(113, 8).Value = 1864.5
$ws.Cells.Item(113, 9).Value = 1877.2222
$ws.Cells.Item(113, 10).Value = 1750
$ws.Cells.Item(113, 11).Value = 1877.2222
$ws.Cells.Item(113, 12).Value = 1750
$ws.Cells.Item(113, 13).Value = 292.7778000000001
$ws.Cells.Item(113, 14).Value = -6090
# row 122
$ws.Cells.Item(122, 8).Value = 3808.2307
$ws.Cells.Item(122, 9).Value = 3808.2307
$ws.Cells.Item(122, 11).Value = 11424.6921
$ws.Cells.Item(122, 13).Value = -8974.6921
# row 132
$ws.Cells.Item(132, 8).Value = 2723.2368
$ws.Cells.Item(132, 9).Value = 2561.5417
$ws.Cells.Item(132, 10).Value = 3000.4285
$ws.Cells.Item(132, 11).Value = 7684.625100000001
$ws.Cells.Item(132, 12).Value = 9001.2855
$ws.Cells.Item(132, 13).Value = -5154.625100000001
$ws.Cells.Item(132, 14).Value = -14061.2855

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Cells.Item(22, 8).Value = 2000
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 2000
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 2000
$ws.Cells.Item(22, 13).Value = $null
$ws.Cells.Item(22, 14).Value = -2590
# row 27
$ws.Cells.Item(27, 8).Value = 2000
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 2000
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 2000
$ws.Cells.Item(27, 13).Value = $null
$ws.Cells.Item(27, 14).Value = -2214
# row 40
$ws.Cells.Item(40, 8).Value = 5180.5
$ws.Cells.Item(40, 9).Value = 2479.5557
$ws.Cells.Item(40, 11).Value = 2479.5557
$ws.Cells.Item(40, 13).Value = -2343.5557
# row 119
$ws.Cells.Item(119, 8).Value = 30000
$ws.Cells.Item(119, 10).Value = 30000
$ws.Cells.Item(119, 12).Value = 30000
$ws.Cells.Item(119, 14).Value = -39676
# row 122
$ws.Cells.Item(122, 8).Value = 16668386
$ws.Cells.Item(122, 9).Value = 22728890
$ws.Cells.Item(122, 11).Value = 68186670
$ws.Cells.Item(122, 13).Value = -68184220
# row 132
$ws.Cells.Item(132, 8).Value = 3476.2354
$ws.Cells.Item(132, 9).Value = 4739.6
$ws.Cells.Item(132, 11).Value = 14218.8
$ws.Cells.Item(132, 13).Value = -11688.8

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 100
$ws.Cells.Item(100, 8).Value = 434.7143
$ws.Cells.Item(100, 9).Value = 438
$ws.Cells.Item(100, 10).Value = 426.5
$ws.Cells.Item(100, 11).Value = 876
$ws.Cells.Item(100, 12).Value = 853
$ws.Cells.Item(100, 13).Value = -335
$ws.Cells.Item(100, 14).Value = -1935
# row 132
$ws.Cells.Item(132, 8).Value = 2045.3889
$ws.Cells.Item(132, 9).Value = 1901.6897
$ws.Cells.Item(132, 10).Value = 2640.7144
$ws.Cells.Item(132, 11).Value = 5705.0691
$ws.Cells.Item(132, 12).Value = 7922.1432
$ws.Cells.Item(132, 13).Value = -3175.0691
$ws.Cells.Item(132, 14).Value = -12982.1432
